# Auto-generated: swap/rotate match-row data blocks (columns B:AD, excluding C/D which are
# identical within each group) per the commit diff for South Korea K3 League.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Cells.Item(14, 2).Value = 6352251
$ws.Cells.Item(14, 5).Value = "Ulsan Citizen FC"
$ws.Cells.Item(14, 6).Value = "Pocheon Citizen FC"
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 1
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = "A"
$ws.Cells.Item(14, 12).Value = 1.909
$ws.Cells.Item(14, 13).Value = 3.25
$ws.Cells.Item(14, 14).Value = 3.5
$ws.Cells.Item(14, 15).Value = 2.25
$ws.Cells.Item(14, 16).Value = 3
$ws.Cells.Item(14, 17).Value = 2.9
$ws.Cells.Item(14, 18).Value = -0.25
$ws.Cells.Item(14, 19).Value = 2.025
$ws.Cells.Item(14, 20).Value = 1.775
$ws.Cells.Item(14, 21).Value = 2
$ws.Cells.Item(14, 22).Value = 1.95
$ws.Cells.Item(14, 23).Value = 1.85
$ws.Cells.Item(14, 24).Value = -1
$ws.Cells.Item(14, 25).Value = -1
$ws.Cells.Item(14, 26).Value = 1.9
$ws.Cells.Item(14, 27).Value = -1
$ws.Cells.Item(14, 28).Value = 0.7749999999999999
$ws.Cells.Item(14, 29).Value = -1
$ws.Cells.Item(14, 30).Value = 0.8500000000000001

# Row 15
$ws.Cells.Item(15, 2).Value = 6352908
$ws.Cells.Item(15, 5).Value = "Gangneung City"
$ws.Cells.Item(15, 6).Value = "Siheung City AC"
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = "D"
$ws.Cells.Item(15, 12).Value = 2.5
$ws.Cells.Item(15, 13).Value = 3.2
$ws.Cells.Item(15, 14).Value = 2.5
$ws.Cells.Item(15, 15).Value = 2.625
$ws.Cells.Item(15, 16).Value = 3
$ws.Cells.Item(15, 17).Value = 2.5
$ws.Cells.Item(15, 18).Value = 0
$ws.Cells.Item(15, 19).Value = 1.95
$ws.Cells.Item(15, 20).Value = 1.85
$ws.Cells.Item(15, 21).Value = 2.25
$ws.Cells.Item(15, 22).Value = 2.025
$ws.Cells.Item(15, 23).Value = 1.775
$ws.Cells.Item(15, 24).Value = -1
$ws.Cells.Item(15, 25).Value = 2
$ws.Cells.Item(15, 26).Value = -1
$ws.Cells.Item(15, 27).Value = 0
$ws.Cells.Item(15, 28).Value = 0
$ws.Cells.Item(15, 29).Value = -1
$ws.Cells.Item(15, 30).Value = 0.7749999999999999

# Row 25
$ws.Cells.Item(25, 2).Value = 6353326
$ws.Cells.Item(25, 5).Value = "Yangpyeong FC"
$ws.Cells.Item(25, 6).Value = "Yangju Citizen"
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 1
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = "A"
$ws.Cells.Item(25, 12).Value = 2.2
$ws.Cells.Item(25, 13).Value = 3.25
$ws.Cells.Item(25, 14).Value = 2.8
$ws.Cells.Item(25, 15).Value = 2.3
$ws.Cells.Item(25, 16).Value = 3.25
$ws.Cells.Item(25, 17).Value = 2.625
$ws.Cells.Item(25, 18).Value = 0
$ws.Cells.Item(25, 19).Value = 1.775
$ws.Cells.Item(25, 20).Value = 2.025
$ws.Cells.Item(25, 21).Value = 2.25
$ws.Cells.Item(25, 22).Value = 1.925
$ws.Cells.Item(25, 23).Value = 1.875
$ws.Cells.Item(25, 24).Value = -1
$ws.Cells.Item(25, 25).Value = -1
$ws.Cells.Item(25, 26).Value = 1.625
$ws.Cells.Item(25, 27).Value = -1
$ws.Cells.Item(25, 28).Value = 1.025
$ws.Cells.Item(25, 29).Value = -1
$ws.Cells.Item(25, 30).Value = 0.875

# Row 26
$ws.Cells.Item(26, 2).Value = 6352256
$ws.Cells.Item(26, 5).Value = "Gyeongju HNP"
$ws.Cells.Item(26, 6).Value = "Pocheon Citizen FC"
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 8).Value = 3
$ws.Cells.Item(26, 9).Value = 1
$ws.Cells.Item(26, 10).Value = 3
$ws.Cells.Item(26, 11).Value = "A"
$ws.Cells.Item(26, 12).Value = 2.625
$ws.Cells.Item(26, 13).Value = 3.2
$ws.Cells.Item(26, 14).Value = 2.375
$ws.Cells.Item(26, 15).Value = 2.1
$ws.Cells.Item(26, 16).Value = 3.2
$ws.Cells.Item(26, 17).Value = 3.1
$ws.Cells.Item(26, 18).Value = -0.25
$ws.Cells.Item(26, 19).Value = 1.875
$ws.Cells.Item(26, 20).Value = 1.925
$ws.Cells.Item(26, 21).Value = 2.25
$ws.Cells.Item(26, 22).Value = 2
$ws.Cells.Item(26, 23).Value = 1.8
$ws.Cells.Item(26, 24).Value = -1
$ws.Cells.Item(26, 25).Value = -1
$ws.Cells.Item(26, 26).Value = 2.1
$ws.Cells.Item(26, 27).Value = -1
$ws.Cells.Item(26, 28).Value = 0.925
$ws.Cells.Item(26, 29).Value = 1
$ws.Cells.Item(26, 30).Value = -1

# Row 53
$ws.Cells.Item(53, 2).Value = 6352927
$ws.Cells.Item(53, 5).Value = "Mokpo City"
$ws.Cells.Item(53, 6).Value = "Siheung City AC"
$ws.Cells.Item(53, 7).Value = 4
$ws.Cells.Item(53, 8).Value = 3
$ws.Cells.Item(53, 9).Value = 2
$ws.Cells.Item(53, 10).Value = 1
$ws.Cells.Item(53, 11).Value = "H"
$ws.Cells.Item(53, 12).Value = 1.909
$ws.Cells.Item(53, 13).Value = 3.2
$ws.Cells.Item(53, 14).Value = 3.6
$ws.Cells.Item(53, 15).Value = 1.85
$ws.Cells.Item(53, 16).Value = 3.5
$ws.Cells.Item(53, 17).Value = 3.6
$ws.Cells.Item(53, 18).Value = -0.5
$ws.Cells.Item(53, 19).Value = 1.9
$ws.Cells.Item(53, 20).Value = 1.9
$ws.Cells.Item(53, 21).Value = 2.5
$ws.Cells.Item(53, 22).Value = 1.975
$ws.Cells.Item(53, 23).Value = 1.825
$ws.Cells.Item(53, 24).Value = 0.8500000000000001
$ws.Cells.Item(53, 25).Value = -1
$ws.Cells.Item(53, 26).Value = -1
$ws.Cells.Item(53, 27).Value = 0.8999999999999999
$ws.Cells.Item(53, 28).Value = -1
$ws.Cells.Item(53, 29).Value = 0.9750000000000001
$ws.Cells.Item(53, 30).Value = -1

# Row 54
$ws.Cells.Item(54, 2).Value = 6352270
$ws.Cells.Item(54, 5).Value = "Paju Citizen FC"
$ws.Cells.Item(54, 6).Value = "Pocheon Citizen FC"
$ws.Cells.Item(54, 7).Value = 3
$ws.Cells.Item(54, 8).Value = 1
$ws.Cells.Item(54, 9).Value = 2
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 11).Value = "H"
$ws.Cells.Item(54, 12).Value = 2
$ws.Cells.Item(54, 13).Value = 3.1
$ws.Cells.Item(54, 14).Value = 3.4
$ws.Cells.Item(54, 15).Value = 1.833
$ws.Cells.Item(54, 16).Value = 3.25
$ws.Cells.Item(54, 17).Value = 4.2
$ws.Cells.Item(54, 18).Value = -0.5
$ws.Cells.Item(54, 19).Value = 1.85
$ws.Cells.Item(54, 20).Value = 1.95
$ws.Cells.Item(54, 21).Value = 2.25
$ws.Cells.Item(54, 22).Value = 2.025
$ws.Cells.Item(54, 23).Value = 1.775
$ws.Cells.Item(54, 24).Value = 0.833
$ws.Cells.Item(54, 25).Value = -1
$ws.Cells.Item(54, 26).Value = -1
$ws.Cells.Item(54, 27).Value = 0.8500000000000001
$ws.Cells.Item(54, 28).Value = -1
$ws.Cells.Item(54, 29).Value = 1.025
$ws.Cells.Item(54, 30).Value = -1

# Row 93
$ws.Cells.Item(93, 2).Value = 6352290
$ws.Cells.Item(93, 5).Value = "Daejeon Korail"
$ws.Cells.Item(93, 6).Value = "Paju Citizen FC"
$ws.Cells.Item(93, 7).Value = 1
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = "H"
$ws.Cells.Item(93, 12).Value = 2.3
$ws.Cells.Item(93, 13).Value = 3
$ws.Cells.Item(93, 14).Value = 2.875
$ws.Cells.Item(93, 15).Value = 2.5
$ws.Cells.Item(93, 16).Value = 3
$ws.Cells.Item(93, 17).Value = 2.625
$ws.Cells.Item(93, 18).Value = 0
$ws.Cells.Item(93, 19).Value = 1.85
$ws.Cells.Item(93, 20).Value = 1.95
$ws.Cells.Item(93, 21).Value = 2
$ws.Cells.Item(93, 22).Value = 1.925
$ws.Cells.Item(93, 23).Value = 1.875
$ws.Cells.Item(93, 24).Value = 1.5
$ws.Cells.Item(93, 25).Value = -1
$ws.Cells.Item(93, 26).Value = -1
$ws.Cells.Item(93, 27).Value = 0.8500000000000001
$ws.Cells.Item(93, 28).Value = -1
$ws.Cells.Item(93, 29).Value = -1
$ws.Cells.Item(93, 30).Value = 0.875

# Row 94
$ws.Cells.Item(94, 2).Value = 6352289
$ws.Cells.Item(94, 5).Value = "Changwon City"
$ws.Cells.Item(94, 6).Value = "Gimhae City"
$ws.Cells.Item(94, 7).Value = 3
$ws.Cells.Item(94, 8).Value = 2
$ws.Cells.Item(94, 9).Value = 3
$ws.Cells.Item(94, 10).Value = 1
$ws.Cells.Item(94, 11).Value = "H"
$ws.Cells.Item(94, 12).Value = 4.2
$ws.Cells.Item(94, 13).Value = 3.6
$ws.Cells.Item(94, 14).Value = 1.666
$ws.Cells.Item(94, 15).Value = 4.5
$ws.Cells.Item(94, 16).Value = 3.8
$ws.Cells.Item(94, 17).Value = 1.571
$ws.Cells.Item(94, 18).Value = 1
$ws.Cells.Item(94, 19).Value = 1.775
$ws.Cells.Item(94, 20).Value = 2.025
$ws.Cells.Item(94, 21).Value = 2.25
$ws.Cells.Item(94, 22).Value = 1.875
$ws.Cells.Item(94, 23).Value = 1.925
$ws.Cells.Item(94, 24).Value = 3.5
$ws.Cells.Item(94, 25).Value = -1
$ws.Cells.Item(94, 26).Value = -1
$ws.Cells.Item(94, 27).Value = 0.7749999999999999
$ws.Cells.Item(94, 28).Value = -1
$ws.Cells.Item(94, 29).Value = 0.875
$ws.Cells.Item(94, 30).Value = -1

# Row 103
$ws.Cells.Item(103, 2).Value = 6352295
$ws.Cells.Item(103, 5).Value = "Gimhae City"
$ws.Cells.Item(103, 6).Value = "Gyeongju HNP"
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 1
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = 1
$ws.Cells.Item(103, 11).Value = "A"
$ws.Cells.Item(103, 12).Value = 1.833
$ws.Cells.Item(103, 13).Value = 3.1
$ws.Cells.Item(103, 14).Value = 4
$ws.Cells.Item(103, 15).Value = 2.25
$ws.Cells.Item(103, 16).Value = 3.1
$ws.Cells.Item(103, 17).Value = 3.1
$ws.Cells.Item(103, 18).Value = -0.25
$ws.Cells.Item(103, 19).Value = 1.975
$ws.Cells.Item(103, 20).Value = 1.825
$ws.Cells.Item(103, 21).Value = 2.25
$ws.Cells.Item(103, 22).Value = 1.95
$ws.Cells.Item(103, 23).Value = 1.85
$ws.Cells.Item(103, 24).Value = -1
$ws.Cells.Item(103, 25).Value = -1
$ws.Cells.Item(103, 26).Value = 2.1
$ws.Cells.Item(103, 27).Value = -1
$ws.Cells.Item(103, 28).Value = 0.825
$ws.Cells.Item(103, 29).Value = -1
$ws.Cells.Item(103, 30).Value = 0.8500000000000001

# Row 104
$ws.Cells.Item(104, 2).Value = 6352294
$ws.Cells.Item(104, 5).Value = "Ulsan Citizen FC"
$ws.Cells.Item(104, 6).Value = "Busan Trans Corp"
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 3
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 11).Value = "A"
$ws.Cells.Item(104, 12).Value = 2.1
$ws.Cells.Item(104, 13).Value = 3
$ws.Cells.Item(104, 14).Value = 3.25
$ws.Cells.Item(104, 15).Value = 2.4
$ws.Cells.Item(104, 16).Value = 2.9
$ws.Cells.Item(104, 17).Value = 2.8
$ws.Cells.Item(104, 18).Value = 0
$ws.Cells.Item(104, 19).Value = 1.775
$ws.Cells.Item(104, 20).Value = 2.025
$ws.Cells.Item(104, 21).Value = 2.25
$ws.Cells.Item(104, 22).Value = 1.95
$ws.Cells.Item(104, 23).Value = 1.85
$ws.Cells.Item(104, 24).Value = -1
$ws.Cells.Item(104, 25).Value = -1
$ws.Cells.Item(104, 26).Value = 1.8
$ws.Cells.Item(104, 27).Value = -1
$ws.Cells.Item(104, 28).Value = 1.025
$ws.Cells.Item(104, 29).Value = 0.95
$ws.Cells.Item(104, 30).Value = -1

# Row 110
$ws.Cells.Item(110, 2).Value = 6352949
$ws.Cells.Item(110, 5).Value = "Mokpo City"
$ws.Cells.Item(110, 6).Value = "Chuncheon FC"
$ws.Cells.Item(110, 7).Value = 3
$ws.Cells.Item(110, 8).Value = 3
$ws.Cells.Item(110, 9).Value = 3
$ws.Cells.Item(110, 10).Value = 1
$ws.Cells.Item(110, 11).Value = "D"
$ws.Cells.Item(110, 12).Value = 1.571
$ws.Cells.Item(110, 13).Value = 3.6
$ws.Cells.Item(110, 14).Value = 5
$ws.Cells.Item(110, 15).Value = 1.363
$ws.Cells.Item(110, 16).Value = 4.75
$ws.Cells.Item(110, 17).Value = 6
$ws.Cells.Item(110, 18).Value = -1.25
$ws.Cells.Item(110, 19).Value = 1.8
$ws.Cells.Item(110, 20).Value = 2
$ws.Cells.Item(110, 21).Value = 2.5
$ws.Cells.Item(110, 22).Value = 1.85
$ws.Cells.Item(110, 23).Value = 1.95
$ws.Cells.Item(110, 24).Value = -1
$ws.Cells.Item(110, 25).Value = 3.75
$ws.Cells.Item(110, 26).Value = -1
$ws.Cells.Item(110, 27).Value = -1
$ws.Cells.Item(110, 28).Value = 1
$ws.Cells.Item(110, 29).Value = 0.8500000000000001
$ws.Cells.Item(110, 30).Value = -1

# Row 111
$ws.Cells.Item(111, 2).Value = 6352948
$ws.Cells.Item(111, 5).Value = "Gyeongju HNP"
$ws.Cells.Item(111, 6).Value = "Siheung City AC"
$ws.Cells.Item(111, 7).Value = 3
$ws.Cells.Item(111, 8).Value = 1
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 1
$ws.Cells.Item(111, 11).Value = "H"
$ws.Cells.Item(111, 12).Value = 2
$ws.Cells.Item(111, 13).Value = 3.1
$ws.Cells.Item(111, 14).Value = 3.4
$ws.Cells.Item(111, 15).Value = 2.15
$ws.Cells.Item(111, 16).Value = 3.1
$ws.Cells.Item(111, 17).Value = 3.1
$ws.Cells.Item(111, 18).Value = -0.25
$ws.Cells.Item(111, 19).Value = 1.925
$ws.Cells.Item(111, 20).Value = 1.875
$ws.Cells.Item(111, 21).Value = 2.25
$ws.Cells.Item(111, 22).Value = 1.9
$ws.Cells.Item(111, 23).Value = 1.9
$ws.Cells.Item(111, 24).Value = 1.15
$ws.Cells.Item(111, 25).Value = -1
$ws.Cells.Item(111, 26).Value = -1
$ws.Cells.Item(111, 27).Value = 0.925
$ws.Cells.Item(111, 28).Value = -1
$ws.Cells.Item(111, 29).Value = 0.8999999999999999
$ws.Cells.Item(111, 30).Value = -1

# Row 118
$ws.Cells.Item(118, 2).Value = 6353333
$ws.Cells.Item(118, 5).Value = "Gyeongju HNP"
$ws.Cells.Item(118, 6).Value = "Busan Trans Corp"
$ws.Cells.Item(118, 7).Value = 3
$ws.Cells.Item(118, 8).Value = 1
$ws.Cells.Item(118, 9).Value = 2
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = "H"
$ws.Cells.Item(118, 12).Value = 2.375
$ws.Cells.Item(118, 13).Value = 2.8
$ws.Cells.Item(118, 14).Value = 3
$ws.Cells.Item(118, 15).Value = 2.15
$ws.Cells.Item(118, 16).Value = 3
$ws.Cells.Item(118, 17).Value = 3.2
$ws.Cells.Item(118, 18).Value = -0.25
$ws.Cells.Item(118, 19).Value = 1.925
$ws.Cells.Item(118, 20).Value = 1.875
$ws.Cells.Item(118, 21).Value = 2.25
$ws.Cells.Item(118, 22).Value = 1.85
$ws.Cells.Item(118, 23).Value = 1.95
$ws.Cells.Item(118, 24).Value = 1.15
$ws.Cells.Item(118, 25).Value = -1
$ws.Cells.Item(118, 26).Value = -1
$ws.Cells.Item(118, 27).Value = 0.925
$ws.Cells.Item(118, 28).Value = -1
$ws.Cells.Item(118, 29).Value = 0.8500000000000001
$ws.Cells.Item(118, 30).Value = -1

# Row 119
$ws.Cells.Item(119, 2).Value = 6352951
$ws.Cells.Item(119, 5).Value = "Paju Citizen FC"
$ws.Cells.Item(119, 6).Value = "Yangju Citizen"
$ws.Cells.Item(119, 7).Value = 4
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 9).Value = 1
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 11).Value = "H"
$ws.Cells.Item(119, 12).Value = 1.8
$ws.Cells.Item(119, 13).Value = 3.4
$ws.Cells.Item(119, 14).Value = 3.8
$ws.Cells.Item(119, 15).Value = 1.65
$ws.Cells.Item(119, 16).Value = 3.6
$ws.Cells.Item(119, 17).Value = 4.333
$ws.Cells.Item(119, 18).Value = -0.75
$ws.Cells.Item(119, 19).Value = 1.825
$ws.Cells.Item(119, 20).Value = 1.975
$ws.Cells.Item(119, 21).Value = 2.5
$ws.Cells.Item(119, 22).Value = 1.875
$ws.Cells.Item(119, 23).Value = 1.925
$ws.Cells.Item(119, 24).Value = 0.6499999999999999
$ws.Cells.Item(119, 25).Value = -1
$ws.Cells.Item(119, 26).Value = -1
$ws.Cells.Item(119, 27).Value = 0.825
$ws.Cells.Item(119, 28).Value = -1
$ws.Cells.Item(119, 29).Value = 0.875
$ws.Cells.Item(119, 30).Value = -1

# Row 124
$ws.Cells.Item(124, 2).Value = 6352952
$ws.Cells.Item(124, 5).Value = "Paju Citizen FC"
$ws.Cells.Item(124, 6).Value = "Mokpo City"
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 1
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 1
$ws.Cells.Item(124, 11).Value = "A"
$ws.Cells.Item(124, 12).Value = 2.6
$ws.Cells.Item(124, 13).Value = 3.2
$ws.Cells.Item(124, 14).Value = 2.4
$ws.Cells.Item(124, 15).Value = 3.1
$ws.Cells.Item(124, 16).Value = 3.3
$ws.Cells.Item(124, 17).Value = 2.05
$ws.Cells.Item(124, 18).Value = 0.25
$ws.Cells.Item(124, 19).Value = 1.95
$ws.Cells.Item(124, 20).Value = 1.85
$ws.Cells.Item(124, 21).Value = 2.25
$ws.Cells.Item(124, 22).Value = 1.8
$ws.Cells.Item(124, 23).Value = 2
$ws.Cells.Item(124, 24).Value = -1
$ws.Cells.Item(124, 25).Value = -1
$ws.Cells.Item(124, 26).Value = 1.05
$ws.Cells.Item(124, 27).Value = -1
$ws.Cells.Item(124, 28).Value = 0.8500000000000001
$ws.Cells.Item(124, 29).Value = -1
$ws.Cells.Item(124, 30).Value = 1

# Row 125
$ws.Cells.Item(125, 2).Value = 6352953
$ws.Cells.Item(125, 5).Value = "Gimhae City"
$ws.Cells.Item(125, 6).Value = "Daejeon Korail"
$ws.Cells.Item(125, 7).Value = 2
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 9).Value = 1
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = "H"
$ws.Cells.Item(125, 12).Value = 2.05
$ws.Cells.Item(125, 13).Value = 3.25
$ws.Cells.Item(125, 14).Value = 3.1
$ws.Cells.Item(125, 15).Value = 1.8
$ws.Cells.Item(125, 16).Value = 3.4
$ws.Cells.Item(125, 17).Value = 3.75
$ws.Cells.Item(125, 18).Value = -0.5
$ws.Cells.Item(125, 19).Value = 1.825
$ws.Cells.Item(125, 20).Value = 1.975
$ws.Cells.Item(125, 21).Value = 2.5
$ws.Cells.Item(125, 22).Value = 2
$ws.Cells.Item(125, 23).Value = 1.8
$ws.Cells.Item(125, 24).Value = 0.8
$ws.Cells.Item(125, 25).Value = -1
$ws.Cells.Item(125, 26).Value = -1
$ws.Cells.Item(125, 27).Value = 0.825
$ws.Cells.Item(125, 28).Value = -1
$ws.Cells.Item(125, 29).Value = -1
$ws.Cells.Item(125, 30).Value = 0.8

# Row 126
$ws.Cells.Item(126, 2).Value = 6352956
$ws.Cells.Item(126, 5).Value = "Ulsan Citizen FC"
$ws.Cells.Item(126, 6).Value = "Yangpyeong FC"
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = "D"
$ws.Cells.Item(126, 12).Value = 2.05
$ws.Cells.Item(126, 13).Value = 3.5
$ws.Cells.Item(126, 14).Value = 2.9
$ws.Cells.Item(126, 15).Value = 2.75
$ws.Cells.Item(126, 16).Value = 3.5
$ws.Cells.Item(126, 17).Value = 2.15
$ws.Cells.Item(126, 18).Value = 0.25
$ws.Cells.Item(126, 19).Value = 1.8
$ws.Cells.Item(126, 20).Value = 2
$ws.Cells.Item(126, 21).Value = 2.5
$ws.Cells.Item(126, 22).Value = 1.95
$ws.Cells.Item(126, 23).Value = 1.85
$ws.Cells.Item(126, 24).Value = -1
$ws.Cells.Item(126, 25).Value = 2.5
$ws.Cells.Item(126, 26).Value = -1
$ws.Cells.Item(126, 27).Value = 0.4
$ws.Cells.Item(126, 28).Value = -0.5
$ws.Cells.Item(126, 29).Value = -1
$ws.Cells.Item(126, 30).Value = 0.8500000000000001

# Row 131
$ws.Cells.Item(131, 2).Value = 7866860
$ws.Cells.Item(131, 5).Value = "Daejeon Korail"
$ws.Cells.Item(131, 6).Value = "Gangneung City"
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = "H"
$ws.Cells.Item(131, 12).Value = 2.1
$ws.Cells.Item(131, 13).Value = 3
$ws.Cells.Item(131, 14).Value = 3.3
$ws.Cells.Item(131, 15).Value = 2.1
$ws.Cells.Item(131, 16).Value = 3
$ws.Cells.Item(131, 17).Value = 3.4
$ws.Cells.Item(131, 18).Value = -0.25
$ws.Cells.Item(131, 19).Value = 1.825
$ws.Cells.Item(131, 20).Value = 1.975
$ws.Cells.Item(131, 21).Value = 2
$ws.Cells.Item(131, 22).Value = 1.85
$ws.Cells.Item(131, 23).Value = 1.95
$ws.Cells.Item(131, 24).Value = 1.1
$ws.Cells.Item(131, 25).Value = -1
$ws.Cells.Item(131, 26).Value = -1
$ws.Cells.Item(131, 27).Value = 0.825
$ws.Cells.Item(131, 28).Value = -1
$ws.Cells.Item(131, 29).Value = -1
$ws.Cells.Item(131, 30).Value = 0.95

# Row 132
$ws.Cells.Item(132, 2).Value = 7873812
$ws.Cells.Item(132, 5).Value = "Hwaseong FC"
$ws.Cells.Item(132, 6).Value = "Daegu FC Reserves"
$ws.Cells.Item(132, 7).Value = 2
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 1
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = "H"
$ws.Cells.Item(132, 12).Value = 1.444
$ws.Cells.Item(132, 13).Value = 4.333
$ws.Cells.Item(132, 14).Value = 5.25
$ws.Cells.Item(132, 15).Value = 1.444
$ws.Cells.Item(132, 16).Value = 4.333
$ws.Cells.Item(132, 17).Value = 5.25
$ws.Cells.Item(132, 18).Value = -1.25
$ws.Cells.Item(132, 19).Value = 1.975
$ws.Cells.Item(132, 20).Value = 1.825
$ws.Cells.Item(132, 21).Value = 2.5
$ws.Cells.Item(132, 22).Value = 1.8
$ws.Cells.Item(132, 23).Value = 2
$ws.Cells.Item(132, 24).Value = 0.444
$ws.Cells.Item(132, 25).Value = -1
$ws.Cells.Item(132, 26).Value = -1
$ws.Cells.Item(132, 27).Value = 0.9750000000000001
$ws.Cells.Item(132, 28).Value = -1
$ws.Cells.Item(132, 29).Value = -1
$ws.Cells.Item(132, 30).Value = 1

# Row 150
$ws.Cells.Item(150, 2).Value = 7867516
$ws.Cells.Item(150, 5).Value = "Yangpyeong FC"
$ws.Cells.Item(150, 6).Value = "Chuncheon FC"
$ws.Cells.Item(150, 7).Value = 1
$ws.Cells.Item(150, 8).Value = 0
$ws.Cells.Item(150, 9).Value = 0
$ws.Cells.Item(150, 10).Value = 0
$ws.Cells.Item(150, 11).Value = "H"
$ws.Cells.Item(150, 12).Value = 2.6
$ws.Cells.Item(150, 13).Value = 3.2
$ws.Cells.Item(150, 14).Value = 2.4
$ws.Cells.Item(150, 15).Value = 2.8
$ws.Cells.Item(150, 16).Value = 3.2
$ws.Cells.Item(150, 17).Value = 2.25
$ws.Cells.Item(150, 18).Value = 0.25
$ws.Cells.Item(150, 19).Value = 1.775
$ws.Cells.Item(150, 20).Value = 2.025
$ws.Cells.Item(150, 21).Value = 2
$ws.Cells.Item(150, 22).Value = 1.85
$ws.Cells.Item(150, 23).Value = 1.95
$ws.Cells.Item(150, 24).Value = 1.8
$ws.Cells.Item(150, 25).Value = -1
$ws.Cells.Item(150, 26).Value = -1
$ws.Cells.Item(150, 27).Value = 0.7749999999999999
$ws.Cells.Item(150, 28).Value = -1
$ws.Cells.Item(150, 29).Value = -1
$ws.Cells.Item(150, 30).Value = 0.95

# Row 151
$ws.Cells.Item(151, 2).Value = 7867518
$ws.Cells.Item(151, 5).Value = "Pocheon Citizen FC"
$ws.Cells.Item(151, 6).Value = "Gyeongju HNP"
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 2
$ws.Cells.Item(151, 9).Value = 0
$ws.Cells.Item(151, 10).Value = 2
$ws.Cells.Item(151, 11).Value = "A"
$ws.Cells.Item(151, 12).Value = 2.4
$ws.Cells.Item(151, 13).Value = 3.2
$ws.Cells.Item(151, 14).Value = 2.6
$ws.Cells.Item(151, 15).Value = 2.9
$ws.Cells.Item(151, 16).Value = 3.2
$ws.Cells.Item(151, 17).Value = 2.15
$ws.Cells.Item(151, 18).Value = 0.25
$ws.Cells.Item(151, 19).Value = 1.85
$ws.Cells.Item(151, 20).Value = 1.95
$ws.Cells.Item(151, 21).Value = 2.25
$ws.Cells.Item(151, 22).Value = 2
$ws.Cells.Item(151, 23).Value = 1.8
$ws.Cells.Item(151, 24).Value = -1
$ws.Cells.Item(151, 25).Value = -1
$ws.Cells.Item(151, 26).Value = 1.15
$ws.Cells.Item(151, 27).Value = -1
$ws.Cells.Item(151, 28).Value = 0.95
$ws.Cells.Item(151, 29).Value = -0.5
$ws.Cells.Item(151, 30).Value = 0.4

# Row 156
$ws.Cells.Item(156, 2).Value = 7867525
$ws.Cells.Item(156, 5).Value = "Siheung City AC"
$ws.Cells.Item(156, 6).Value = "Mokpo City"
$ws.Cells.Item(156, 7).Value = 5
$ws.Cells.Item(156, 8).Value = 0
$ws.Cells.Item(156, 9).Value = 2
$ws.Cells.Item(156, 10).Value = 0
$ws.Cells.Item(156, 11).Value = "H"
$ws.Cells.Item(156, 12).Value = 2.4
$ws.Cells.Item(156, 13).Value = 3.6
$ws.Cells.Item(156, 14).Value = 2.4
$ws.Cells.Item(156, 15).Value = 2.45
$ws.Cells.Item(156, 16).Value = 3.75
$ws.Cells.Item(156, 17).Value = 2.3
$ws.Cells.Item(156, 18).Value = 0
$ws.Cells.Item(156, 19).Value = 1.975
$ws.Cells.Item(156, 20).Value = 1.825
$ws.Cells.Item(156, 21).Value = 2.25
$ws.Cells.Item(156, 22).Value = 1.9
$ws.Cells.Item(156, 23).Value = 1.9
$ws.Cells.Item(156, 24).Value = 1.45
$ws.Cells.Item(156, 25).Value = -1
$ws.Cells.Item(156, 26).Value = -1
$ws.Cells.Item(156, 27).Value = 0.9750000000000001
$ws.Cells.Item(156, 28).Value = -1
$ws.Cells.Item(156, 29).Value = 0.8999999999999999
$ws.Cells.Item(156, 30).Value = -1

# Row 158
$ws.Cells.Item(158, 2).Value = 7867524
$ws.Cells.Item(158, 5).Value = "Hwaseong FC"
$ws.Cells.Item(158, 6).Value = "Paju Citizen FC"
$ws.Cells.Item(158, 7).Value = 2
$ws.Cells.Item(158, 8).Value = 0
$ws.Cells.Item(158, 9).Value = 1
$ws.Cells.Item(158, 10).Value = 0
$ws.Cells.Item(158, 11).Value = "H"
$ws.Cells.Item(158, 12).Value = 1.65
$ws.Cells.Item(158, 13).Value = 3.5
$ws.Cells.Item(158, 14).Value = 5
$ws.Cells.Item(158, 15).Value = 1.6
$ws.Cells.Item(158, 16).Value = 3.6
$ws.Cells.Item(158, 17).Value = 5.5
$ws.Cells.Item(158, 18).Value = -0.75
$ws.Cells.Item(158, 19).Value = 1.775
$ws.Cells.Item(158, 20).Value = 2.025
$ws.Cells.Item(158, 21).Value = 2.25
$ws.Cells.Item(158, 22).Value = 1.825
$ws.Cells.Item(158, 23).Value = 1.975
$ws.Cells.Item(158, 24).Value = 0.6000000000000001
$ws.Cells.Item(158, 25).Value = -1
$ws.Cells.Item(158, 26).Value = -1
$ws.Cells.Item(158, 27).Value = 0.7749999999999999
$ws.Cells.Item(158, 28).Value = -1
$ws.Cells.Item(158, 29).Value = -0.5
$ws.Cells.Item(158, 30).Value = 0.4875

# Row 171
$ws.Cells.Item(171, 2).Value = 7867541
$ws.Cells.Item(171, 5).Value = "Pocheon Citizen FC"
$ws.Cells.Item(171, 6).Value = "Gimhae City"
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 0
$ws.Cells.Item(171, 9).Value = 0
$ws.Cells.Item(171, 10).Value = 0
$ws.Cells.Item(171, 11).Value = "D"
$ws.Cells.Item(171, 12).Value = 2.5
$ws.Cells.Item(171, 13).Value = 3.4
$ws.Cells.Item(171, 14).Value = 2.375
$ws.Cells.Item(171, 15).Value = 3.5
$ws.Cells.Item(171, 16).Value = 3.3
$ws.Cells.Item(171, 17).Value = 1.85
$ws.Cells.Item(171, 18).Value = 0.5
$ws.Cells.Item(171, 19).Value = 1.85
$ws.Cells.Item(171, 20).Value = 1.95
$ws.Cells.Item(171, 21).Value = 2.25
$ws.Cells.Item(171, 22).Value = 1.825
$ws.Cells.Item(171, 23).Value = 1.975
$ws.Cells.Item(171, 24).Value = -1
$ws.Cells.Item(171, 25).Value = 2.3
$ws.Cells.Item(171, 26).Value = -1
$ws.Cells.Item(171, 27).Value = 0.8500000000000001
$ws.Cells.Item(171, 28).Value = -1
$ws.Cells.Item(171, 29).Value = -1
$ws.Cells.Item(171, 30).Value = 0.9750000000000001

# Row 172
$ws.Cells.Item(172, 2).Value = 7867542
$ws.Cells.Item(172, 5).Value = "Busan Trans Corp"
$ws.Cells.Item(172, 6).Value = "Daejeon Korail"
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 0
$ws.Cells.Item(172, 9).Value = 0
$ws.Cells.Item(172, 10).Value = 0
$ws.Cells.Item(172, 11).Value = "D"
$ws.Cells.Item(172, 12).Value = 3.5
$ws.Cells.Item(172, 13).Value = 3.3
$ws.Cells.Item(172, 14).Value = 1.909
$ws.Cells.Item(172, 15).Value = 3.5
$ws.Cells.Item(172, 16).Value = 3.3
$ws.Cells.Item(172, 17).Value = 1.909
$ws.Cells.Item(172, 18).Value = 0.5
$ws.Cells.Item(172, 19).Value = 1.85
$ws.Cells.Item(172, 20).Value = 1.95
$ws.Cells.Item(172, 21).Value = 2.25
$ws.Cells.Item(172, 22).Value = 2
$ws.Cells.Item(172, 23).Value = 1.8
$ws.Cells.Item(172, 24).Value = -1
$ws.Cells.Item(172, 25).Value = 2.3
$ws.Cells.Item(172, 26).Value = -1
$ws.Cells.Item(172, 27).Value = 0.8500000000000001
$ws.Cells.Item(172, 28).Value = -1
$ws.Cells.Item(172, 29).Value = -1
$ws.Cells.Item(172, 30).Value = 0.8

# Row 173
$ws.Cells.Item(173, 2).Value = 7873806
$ws.Cells.Item(173, 5).Value = "Daegu FC Reserves"
$ws.Cells.Item(173, 6).Value = "Paju Citizen FC"
$ws.Cells.Item(173, 7).Value = 2
$ws.Cells.Item(173, 8).Value = 4
$ws.Cells.Item(173, 9).Value = 1
$ws.Cells.Item(173, 10).Value = 4
$ws.Cells.Item(173, 11).Value = "A"
$ws.Cells.Item(173, 12).Value = 3
$ws.Cells.Item(173, 13).Value = 3.25
$ws.Cells.Item(173, 14).Value = 2.1
$ws.Cells.Item(173, 15).Value = 3
$ws.Cells.Item(173, 16).Value = 3.25
$ws.Cells.Item(173, 17).Value = 2.1
$ws.Cells.Item(173, 18).Value = 0.25
$ws.Cells.Item(173, 19).Value = 1.9
$ws.Cells.Item(173, 20).Value = 1.9
$ws.Cells.Item(173, 21).Value = 2.25
$ws.Cells.Item(173, 22).Value = 1.9
$ws.Cells.Item(173, 23).Value = 1.9
$ws.Cells.Item(173, 24).Value = -1
$ws.Cells.Item(173, 25).Value = -1
$ws.Cells.Item(173, 26).Value = 1.1
$ws.Cells.Item(173, 27).Value = -1
$ws.Cells.Item(173, 28).Value = 0.8999999999999999
$ws.Cells.Item(173, 29).Value = 0.8999999999999999
$ws.Cells.Item(173, 30).Value = -1

# Row 196
$ws.Cells.Item(196, 2).Value = 7867566
$ws.Cells.Item(196, 5).Value = "Paju Citizen FC"
$ws.Cells.Item(196, 6).Value = "Yangpyeong FC"
$ws.Cells.Item(196, 7).Value = 1
$ws.Cells.Item(196, 8).Value = 1
$ws.Cells.Item(196, 9).Value = 0
$ws.Cells.Item(196, 10).Value = 0
$ws.Cells.Item(196, 11).Value = "D"
$ws.Cells.Item(196, 12).Value = 1.65
$ws.Cells.Item(196, 13).Value = 3.4
$ws.Cells.Item(196, 14).Value = 4.75
$ws.Cells.Item(196, 15).Value = 1.8
$ws.Cells.Item(196, 16).Value = 3.4
$ws.Cells.Item(196, 17).Value = 3.7
$ws.Cells.Item(196, 18).Value = -0.5
$ws.Cells.Item(196, 19).Value = 1.85
$ws.Cells.Item(196, 20).Value = 1.95
$ws.Cells.Item(196, 21).Value = 2.25
$ws.Cells.Item(196, 22).Value = 1.975
$ws.Cells.Item(196, 23).Value = 1.825
$ws.Cells.Item(196, 24).Value = -1
$ws.Cells.Item(196, 25).Value = 2.4
$ws.Cells.Item(196, 26).Value = -1
$ws.Cells.Item(196, 27).Value = -1
$ws.Cells.Item(196, 28).Value = 0.95
$ws.Cells.Item(196, 29).Value = -0.5
$ws.Cells.Item(196, 30).Value = 0.4125

# Row 198
$ws.Cells.Item(198, 2).Value = 7867564
$ws.Cells.Item(198, 5).Value = "Hwaseong FC"
$ws.Cells.Item(198, 6).Value = "Ulsan Citizen FC"
$ws.Cells.Item(198, 7).Value = 1
$ws.Cells.Item(198, 8).Value = 0
$ws.Cells.Item(198, 9).Value = 0
$ws.Cells.Item(198, 10).Value = 0
$ws.Cells.Item(198, 11).Value = "H"
$ws.Cells.Item(198, 12).Value = 1.75
$ws.Cells.Item(198, 13).Value = 3.3
$ws.Cells.Item(198, 14).Value = 4.2
$ws.Cells.Item(198, 15).Value = 1.727
$ws.Cells.Item(198, 16).Value = 3.4
$ws.Cells.Item(198, 17).Value = 4.2
$ws.Cells.Item(198, 18).Value = -0.75
$ws.Cells.Item(198, 19).Value = 1.95
$ws.Cells.Item(198, 20).Value = 1.85
$ws.Cells.Item(198, 21).Value = 2.25
$ws.Cells.Item(198, 22).Value = 2.025
$ws.Cells.Item(198, 23).Value = 1.775
$ws.Cells.Item(198, 24).Value = 0.7270000000000001
$ws.Cells.Item(198, 25).Value = -1
$ws.Cells.Item(198, 26).Value = -1
$ws.Cells.Item(198, 27).Value = 0.475
$ws.Cells.Item(198, 28).Value = -0.5
$ws.Cells.Item(198, 29).Value = -1
$ws.Cells.Item(198, 30).Value = 0.7749999999999999
